# Ward49.xlsx edit: title-case names/addresses/loc_type_sec values that were
# previously ALL-CAPS, normalize "PLAYGROUND PARK" -> "Playground", and append a
# new Fire Station record (E102) as row 69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 'Sullivan Hs'
$ws.Range("B3").Value = '6631 N Bosworth Ave'
$ws.Range("A4").Value = 'Gale'
$ws.Range("B4").Value = '1631 W Jonquil Ter'
$ws.Range("A5").Value = 'Acero - De La Cruz'
$ws.Range("B5").Value = '7416 N Ridge Blvd'
$ws.Range("A6").Value = 'Field'
$ws.Range("B6").Value = '7019 N Ashland Ave'
$ws.Range("A7").Value = 'Jordan'
$ws.Range("B7").Value = '7414 N Wolcott Ave'
$ws.Range("A8").Value = 'Chicago Math & Science Hs'
$ws.Range("B8").Value = '7212 N Clark St'
$ws.Range("A9").Value = 'Kilmer'
$ws.Range("B9").Value = '6700 N Greenview Ave'
$ws.Range("A10").Value = 'New Field'
$ws.Range("B10").Value = '1707 W Morse Ave'
$ws.Range("A11").Value = 'Hartigan (David) Beach'
$ws.Range("L11").Value = 'Playground'
$ws.Range("A12").Value = 'Howard (Ure) Beach'
$ws.Range("L12").Value = 'Beach'
$ws.Range("A13").Value = 'Howard (Ure) Beach'
$ws.Range("L13").Value = 'Playground'
$ws.Range("A14").Value = 'Griffin (Marion Mahony) Beach'
$ws.Range("L14").Value = 'Beach'
$ws.Range("A15").Value = 'Langdon (Mary Margaret)'
$ws.Range("L15").Value = 'Playground'
$ws.Range("A16").Value = 'Lazarus (Emma)'
$ws.Range("L16").Value = 'Playground'
$ws.Range("A17").Value = 'Leone (Sam) Beach'
$ws.Range("L17").Value = 'Beach'
$ws.Range("A18").Value = 'Leone (Sam) Beach'
$ws.Range("L18").Value = 'Gymnasium'
$ws.Range("A19").Value = 'Leone (Sam) Beach'
$ws.Range("L19").Value = 'Playground'
$ws.Range("A20").Value = 'Loyola'
$ws.Range("L20").Value = 'Basketball Court'
$ws.Range("A21").Value = 'Loyola'
$ws.Range("L21").Value = 'Beach'
$ws.Range("A22").Value = 'Loyola'
$ws.Range("L22").Value = 'Football/Soccer Combo Fld'
$ws.Range("A23").Value = 'Loyola'
$ws.Range("L23").Value = 'Gymnasium'
$ws.Range("A24").Value = 'Loyola'
$ws.Range("L24").Value = 'Baseball Jr/Softball'
$ws.Range("A25").Value = 'Loyola'
$ws.Range("L25").Value = 'Playground'
$ws.Range("A26").Value = 'Loyola'
$ws.Range("L26").Value = 'Baseball Sr'
$ws.Range("A27").Value = 'Loyola'
$ws.Range("L27").Value = 'Tennis Court'
$ws.Range("A28").Value = 'Leone (Sam) Beach'
$ws.Range("L28").Value = 'Boat Launch Non-Motorized'
$ws.Range("A29").Value = 'Columbia Beach'
$ws.Range("L29").Value = 'Beach'
$ws.Range("A30").Value = 'Matanky (Eugene)'
$ws.Range("L30").Value = 'Playground'
$ws.Range("A31").Value = 'Fargo (James) Beach'
$ws.Range("L31").Value = 'Beach'
$ws.Range("A32").Value = 'Goldberg (Louis)'
$ws.Range("L32").Value = 'Playground'
$ws.Range("A33").Value = 'Loyola'
$ws.Range("L33").Value = 'Shuffleboard'
$ws.Range("A34").Value = 'Loyola'
$ws.Range("L34").Value = 'Boxing Center'
$ws.Range("A35").Value = 'Paschen (Christian)'
$ws.Range("L35").Value = 'Volleyball'
$ws.Range("A36").Value = 'White (Willye B.)'
$ws.Range("L36").Value = 'Gymnasium'
$ws.Range("A37").Value = 'Pottawattomie'
$ws.Range("L37").Value = 'Handball/Racquet (In)'
$ws.Range("A38").Value = 'Pottawattomie'
$ws.Range("L38").Value = 'Dog Friendly Area'
$ws.Range("A39").Value = 'Pottawattomie'
$ws.Range("L39").Value = 'Artificial Turf Field'
$ws.Range("A40").Value = 'Loyola'
$ws.Range("L40").Value = 'Nature/Bird Sanctuary'
$ws.Range("A41").Value = 'Juneway Terr. Beach'
$ws.Range("L41").Value = 'Beach'
$ws.Range("A42").Value = 'Prinz (Tobey) Beach'
$ws.Range("L42").Value = 'Beach'
$ws.Range("A43").Value = 'North Shore Beach'
$ws.Range("L43").Value = 'Beach'
$ws.Range("A44").Value = 'Paschen (Christian)'
$ws.Range("L44").Value = 'Basketball Court'
$ws.Range("A45").Value = 'Paschen (Christian)'
$ws.Range("L45").Value = 'Playground'
$ws.Range("A46").Value = 'Pottawattomie'
$ws.Range("L46").Value = 'Basketball Court'
$ws.Range("A47").Value = 'Pottawattomie'
$ws.Range("L47").Value = 'Football/Soccer Combo Fld'
$ws.Range("A48").Value = 'Pottawattomie'
$ws.Range("L48").Value = 'Fitness Center'
$ws.Range("A49").Value = 'Pottawattomie'
$ws.Range("L49").Value = 'Gymnasium'
$ws.Range("A50").Value = 'Pottawattomie'
$ws.Range("L50").Value = 'Baseball Jr/Softball'
$ws.Range("A51").Value = 'Pottawattomie'
$ws.Range("L51").Value = 'Playground'
$ws.Range("A52").Value = 'Pottawattomie'
$ws.Range("L52").Value = 'Spray Feature'
$ws.Range("A53").Value = 'Pottawattomie'
$ws.Range("L53").Value = 'Baseball Sr'
$ws.Range("A54").Value = 'Rogers (Phillip) Beach'
$ws.Range("L54").Value = 'Beach'
$ws.Range("A55").Value = 'Rogers (Phillip) Beach'
$ws.Range("L55").Value = 'Tennis Court'
$ws.Range("A56").Value = 'Touhy (Patrick)'
$ws.Range("L56").Value = 'Football/Soccer Combo Fld'
$ws.Range("A57").Value = 'Touhy (Patrick)'
$ws.Range("L57").Value = 'Baseball Jr/Softball'
$ws.Range("A58").Value = 'Touhy (Patrick)'
$ws.Range("L58").Value = 'Playground'
$ws.Range("A59").Value = 'Touhy (Patrick)'
$ws.Range("L59").Value = 'Tennis Court'
$ws.Range("A60").Value = 'Washington (Harold) Mem.'
$ws.Range("L60").Value = 'Playground'
$ws.Range("A61").Value = 'White (Willye B.)'
$ws.Range("L61").Value = 'Basketball Court'
$ws.Range("A62").Value = 'White (Willye B.)'
$ws.Range("L62").Value = 'Baseball Jr/Softball'
$ws.Range("A63").Value = 'White (Willye B.)'
$ws.Range("L63").Value = 'Playground'
$ws.Range("A64").Value = 'White (Willye B.)'
$ws.Range("L64").Value = 'Fitness Center'
$ws.Range("A65").Value = 'Triangle'
$ws.Range("L65").Value = 'Community Garden'
$ws.Range("A66").Value = 'Dubkin (Leonard)'
$ws.Range("L66").Value = 'Community Garden'
$ws.Range("A67").Value = 'Goldberg (Louis)'
$ws.Range("L67").Value = 'Community Garden'
$ws.Range("A68").Value = 'Loyola'
$ws.Range("L68").Value = 'Community Garden'
$ws.Range("A69").Value = 'E102'
$ws.Range("B69").Value = '7340 N Clark St'
$ws.Range("C69").Value = 'Chicago'
$ws.Range("D69").Value = 'IL'
$ws.Range("E69").Value = '60626'
$ws.Range("H69").Value = 42.01464016068629
$ws.Range("I69").Value = -87.67502361090746
$ws.Range("K69").Value = 'Fire Station'
$ws.Range("L69").Value = 'Fire Station'

Write-Host "Done applying Ward49 fixes."
